# Fix mistake in trial names
#
# Bug: the W-column "sex/trial-type" classifier formula used on sheet
# "Trials_KFJ_base" mapped the (B=1, E=2) case to "F2" instead of "M2".
# That wrong "F2" label then propagated into every X-column
# CONCATENATE(W,"-",U,"-",V) trial name, e.g. "F2-T1-R1" should have
# read "M2-T1-R1".
#
# The worksheet stores this formula as 7 master cells (one plain formula
# in W2, plus six shared-formula groups rooted at W3, W67, W131, W195,
# W259 and W323). Re-assigning `.Formula` across each full shared range
# keeps the same shared-formula grouping/`si` ids Excel already used,
# and lets Excel recalculate every dependent cached value (W and X
# columns) for us.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trials_KFJ_base")

$ranges = @(
    @{ Addr = "W2";       Row = 2   },
    @{ Addr = "W3:W66";   Row = 3   },
    @{ Addr = "W67:W130"; Row = 67  },
    @{ Addr = "W131:W194";Row = 131 },
    @{ Addr = "W195:W258";Row = 195 },
    @{ Addr = "W259:W322";Row = 259 },
    @{ Addr = "W323:W369";Row = 323 }
)

foreach ($r in $ranges) {
    $row = $r.Row
    $formula = '=IF(AND(B' + $row + '=0,E' + $row + '=1),"F1",IF(AND(B' + $row + '=0,E' + $row + '=2),"F2",IF(AND(B' + $row + '=1,E' + $row + '=1),"M1",IF(AND(B' + $row + '=1,E' + $row + '=2),"M2","?"))))'
    $ws.Range($r.Addr).Formula = $formula
}

# Match the author's recorded view/selection state on the sheet: the
# active cell moved from E18 to V9 (still in the frozen bottom-left
# pane), and the sheet's scroll position shifted one column right.
$ws.Activate()
$ws.Range("V9").Select()

$window = $excel.ActiveWindow
$window.ScrollColumn = 5
